# Refresh the cryptos list (prices + 1h volume %) as per the
# "Updated cryptos list ... with GitHub Actions" automated commit.
# Note: Price column (D) cells are stored as plain text (e.g. "60.785.96",
# "1.00"), so numeric-looking values are written with a leading apostrophe
# to force text interpretation, then the cell style is reset to "Normal"
# so no stray quote-prefix style/number format is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'60.785.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "'3.366.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'568.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").Value = "'7.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.23%  "

$ws.Range("E11").Value = "  -4.17%  "

$ws.Range("D12").Value = "'3.942.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("D14").Value = "'27.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").Value = "'3.358.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("D17").Value = "'60.925.83"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "  -2.49%  "

$ws.Range("E19").Value = "  -3.72%  "

$ws.Range("D20").Value = "'8.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("D21").Value = "'382.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("E22").Value = "  +2.81%  "

$ws.Range("D23").Value = "'0.548"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.191"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.61%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000109"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.63%  "

$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("E28").Value = "  -3.90%  "

$ws.Range("D29").Value = "'7.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.42%  "

$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("E32").Value = "  -6.27%  "

$ws.Range("D33").Value = "'22.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.29%  "

$ws.Range("D34").Value = "'167.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("D35").Value = "'6.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("E36").Value = "  -2.24%  "

$ws.Range("D37").Value = "'3.401.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("E38").Value = "  -3.36%  "

$ws.Range("E39").Value = "  -2.91%  "

$ws.Range("D40").Value = "'25.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.60%  "

$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("D42").Value = "'4.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.99%  "

$ws.Range("E43").Value = "  -3.60%  "

$ws.Range("D44").Value = "'2.456.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.49%  "

$ws.Range("E45").Value = "  -2.64%  "

$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "'6.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.35%  "

$ws.Range("D48").Value = "'22.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.81%  "

$ws.Range("E49").Value = "  -5.07%  "

$ws.Range("E50").Value = "  -5.02%  "

$ws.Range("E51").Value = "  -3.26%  "
